# Auto-generated edit script: updates the cryptos list Price (D) and Volume(1h) (E)
# columns for rows 2-51 to reflect the latest scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.764.58"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.313.48"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.29"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("E12").Value = "  -3.35%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "2.653.10"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.876"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.62%  "
$ws.Range("D17").Value = "2.324.13"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "43.717.80"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.38"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.80"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.48"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.39"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.85"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -8.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.31"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0906"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.47"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("E36").Value = "  +2.79%  "
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.246"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.12%  "
$ws.Range("E41").Value = "  +8.36%  "
$ws.Range("E42").Value = "  +18.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.18"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.37"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("E47").Value = "  +3.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.20"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("E50").Value = "  +15.15%  "
$ws.Range("D51").Value = "2.539.43"
$ws.Range("E51").Value = "  +3.95%  "
